$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume update (GitHub Actions daily refresh)
# D-column cells whose text would otherwise be mis-parsed as numbers by
# the COM Value setter are pre-formatted as Text ("@") so values such as
# "1.000" / "0.06560" keep their exact digits instead of being coerced to
# doubles (which would silently drop the significant trailing zeros).

# Row 2
$ws.Range("D2").Value = "27.316.48"
$ws.Range("E2").Value = "  +6.51%  "

# Row 3
$ws.Range("D3").Value = "1.811.61"
$ws.Range("E3").Value = "  +6.22%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.67%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "341.92"
$ws.Range("E5").Value = "  +3.49%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3856"
$ws.Range("E7").Value = "  +5.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.15"
$ws.Range("E8").Value = "  +3.71%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3524"
$ws.Range("E9").Value = "  +7.25%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.234"
$ws.Range("E10").Value = "  +6.03%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07768"
$ws.Range("E11").Value = "  +6.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  +0.79%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.59"
$ws.Range("E13").Value = "  +13.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.637"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.237"
$ws.Range("E15").Value = "  +6.59%  "

# Row 16
$ws.Range("D16").Value = "1.812.74"
$ws.Range("E16").Value = "  +6.54%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001126"
$ws.Range("E17").Value = "  +5.57%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06766"
$ws.Range("E18").Value = "  +2.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.11"
$ws.Range("E19").Value = "  +7.64%  "

# Row 20
$ws.Range("E20").Value = "  +0.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.96"
$ws.Range("E21").Value = "  +11.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.608"
$ws.Range("E22").Value = "  +9.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.10"
$ws.Range("E23").Value = "  +0.53%  "

# Row 24
$ws.Range("D24").Value = "27.352.73"
$ws.Range("E24").Value = "  +6.68%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.476"
$ws.Range("E25").Value = "  +0.87%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.708"
$ws.Range("E26").Value = "  +9.57%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.10"
$ws.Range("E27").Value = "  +15.72%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.511"
$ws.Range("E28").Value = "  +19.00%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.38"
$ws.Range("E29").Value = "  +2.42%  "

# Row 30
$ws.Range("D30").Value = "2.017.39"
$ws.Range("E30").Value = "  +6.64%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "137.52"
$ws.Range("E31").Value = "  +7.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.418"
$ws.Range("E32").Value = "  +7.74%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.150"
$ws.Range("E33").Value = "  +1.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.90"
$ws.Range("E34").Value = "  +9.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08847"
$ws.Range("E35").Value = "  +4.53%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.727"
$ws.Range("E36").Value = "  +2.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.674"
$ws.Range("E37").Value = "  +7.22%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06560"
$ws.Range("E38").Value = "  +5.89%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02427"
$ws.Range("E39").Value = "  +7.50%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2270"
$ws.Range("E40").Value = "  +7.30%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6916"
$ws.Range("E41").Value = "  +13.89%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.007"
$ws.Range("E42").Value = "  +6.48%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.256"
$ws.Range("E43").Value = "  -0.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.10"
$ws.Range("E44").Value = "  +7.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6501"
$ws.Range("E45").Value = "  +11.42%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.018"
$ws.Range("E47").Value = "  +4.73%  "

# Row 48
$ws.Range("E48").Value = "  +9.28%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "133.45"
$ws.Range("E49").Value = "  +6.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07347"
$ws.Range("E50").Value = "  +1.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.88"
$ws.Range("E51").Value = "  +6.28%  "
